# Add two new weekly price records at the top of the data table (rows 341-342),
# pushing the existing rows 341-357 down to 343-359.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 341.
$ws.Rows("341:342").Insert()

# --- New row 341 ---
$ws.Cells.Item(341, 1).Value  = 7
$ws.Cells.Item(341, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(341, 3).Value  = "Ñuble"
$ws.Cells.Item(341, 4).Value  = 45147
$ws.Cells.Item(341, 5).Value  = 16
$ws.Cells.Item(341, 6).Value  = 100112045
$ws.Cells.Item(341, 7).Value  = "Zapallo"
$ws.Cells.Item(341, 8).Value  = "Camote"
$ws.Cells.Item(341, 9).Value  = "1a (guarda)"
$ws.Cells.Item(341, 10).Value = 300
$ws.Cells.Item(341, 11).Value = 650
$ws.Cells.Item(341, 12).Value = 650
$ws.Cells.Item(341, 13).Value = 650
$ws.Cells.Item(341, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(341, 15).Value = "Región del Maule"
$ws.Cells.Item(341, 16).Value = 650
$ws.Cells.Item(341, 17).Value = 1
$ws.Cells.Item(341, 18).Value = "Hortaliza"

# --- New row 342 ---
$ws.Cells.Item(342, 1).Value  = 7
$ws.Cells.Item(342, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(342, 3).Value  = "Ñuble"
$ws.Cells.Item(342, 4).Value  = 45147
$ws.Cells.Item(342, 5).Value  = 16
$ws.Cells.Item(342, 6).Value  = 100112045
$ws.Cells.Item(342, 7).Value  = "Zapallo"
$ws.Cells.Item(342, 8).Value  = "Paine"
$ws.Cells.Item(342, 9).Value  = "1a (guarda)"
$ws.Cells.Item(342, 10).Value = 200
$ws.Cells.Item(342, 11).Value = 350
$ws.Cells.Item(342, 12).Value = 350
$ws.Cells.Item(342, 13).Value = 350
$ws.Cells.Item(342, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(342, 15).Value = "Región del Maule"
$ws.Cells.Item(342, 16).Value = 350
$ws.Cells.Item(342, 17).Value = 1
$ws.Cells.Item(342, 18).Value = "Hortaliza"

# Make sure the date cells keep the date style/format used by the rest of column D.
$ws.Range("D341").NumberFormat = $ws.Range("D343").NumberFormat
$ws.Range("D342").NumberFormat = $ws.Range("D343").NumberFormat
